$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Country name re-sort (sharedStrings reorder) ---
$ws.Range("A90").Value = "Libano"
$ws.Range("A91").Value = "Malasia"
$ws.Range("A146").Value = "Republica de Chipre"
$ws.Range("A147").Value = "Georgia"
$ws.Range("A150").Value = "Botsuana"
$ws.Range("A151").Value = "Burkina Faso"
$ws.Range("A152").Value = "Liberia"
$ws.Range("A213").Value = "Islas Malvinas"
$ws.Range("A214").Value = "Montserrat"

# --- Updated "Datos actualizados" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 17 de Agosto de 2020 a las 19:06"

# --- Updated COVID case numbers (columns B-H) ---
$ws.Range("B4").Value = 5575688
$ws.Range("C4").Value = 9056
$ws.Range("D4").Value = 2925895
$ws.Range("E4").Value = 2476576
$ws.Range("G4").Value = 89
$ws.Range("H4").Value = 173217
$ws.Range("B5").Value = 3343925
$ws.Range("C5").Value = 3728
$ws.Range("E5").Value = 803415
$ws.Range("G5").Value = 175
$ws.Range("H5").Value = 108054
$ws.Range("B6").Value = 2690831
$ws.Range("C6").Value = 43515
$ws.Range("D6").Value = 1965798
$ws.Range("E6").Value = 673238
$ws.Range("G6").Value = 750
$ws.Range("H6").Value = 51795
$ws.Range("D17").Value = 217850
$ws.Range("E17").Value = 70969
$ws.Range("G17").Value = 47
$ws.Range("H17").Value = 5750
$ws.Range("B22").Value = 226351
$ws.Range("C22").Value = 1354
$ws.Range("E22").Value = 14158
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = 9293
$ws.Range("B23").Value = 219029
$ws.Range("C23").Value = 493
$ws.Range("E23").Value = 104752
$ws.Range("G23").Value = 19
$ws.Range("H23").Value = 30429
$ws.Range("B33").Value = 94277
$ws.Range("C33").Value = 1597
$ws.Range("D33").Value = 70267
$ws.Range("E33").Value = 23318
$ws.Range("G33").Value = 7
$ws.Range("H33").Value = 692
$ws.Range("D48").Value = 52350
$ws.Range("E48").Value = 3461
$ws.Range("E53").Value = 3534
$ws.Range("G53").Value = 3
$ws.Range("H53").Value = 173
$ws.Range("B58").Value = 39025
$ws.Range("C58").Value = 442
$ws.Range("D58").Value = 27347
$ws.Range("E58").Value = 10299
$ws.Range("G58").Value = 9
$ws.Range("H58").Value = 1379
$ws.Range("B61").Value = 35702
$ws.Range("C61").Value = 373
$ws.Range("E61").Value = 4493
$ws.Range("G61").Value = 4
$ws.Range("H61").Value = 236
$ws.Range("B73").Value = 23193
$ws.Range("C73").Value = 281
$ws.Range("D73").Value = 10977
$ws.Range("E73").Value = 11598
$ws.Range("B74").Value = 20098
$ws.Range("C74").Value = 86
$ws.Range("D74").Value = 14523
$ws.Range("E74").Value = 5176
$ws.Range("G74").Value = 2
$ws.Range("H74").Value = 399
$ws.Range("B90").Value = 9337
$ws.Range("C90").Value = 456
$ws.Range("D90").Value = 2809
$ws.Range("E90").Value = 6423
$ws.Range("G90").Value = 2
$ws.Range("H90").Value = 105
$ws.Range("B91").Value = 9212
$ws.Range("C91").Value = 12
$ws.Range("D91").Value = 8876
$ws.Range("E91").Value = 211
$ws.Range("H91").Value = 125
$ws.Range("B92").Value = 8622
$ws.Range("C92").Value = 34
$ws.Range("D92").Value = 7985
$ws.Range("E92").Value = 584
$ws.Range("B118").Value = 3364
$ws.Range("C118").Value = 48
$ws.Range("D118").Value = 2692
$ws.Range("E118").Value = 584
$ws.Range("B119").Value = 3257
$ws.Range("C119").Value = 1
$ws.Range("E119").Value = 790
$ws.Range("B124").Value = 2900
$ws.Range("C124").Value = 7
$ws.Range("E124").Value = 213
$ws.Range("B144").Value = 1398
$ws.Range("C144").Value = 20
$ws.Range("D144").Value = 1241
$ws.Range("E144").Value = 146
$ws.Range("B146").Value = 1351
$ws.Range("C146").Value = 12
$ws.Range("D146").Value = 870
$ws.Range("E146").Value = 461
$ws.Range("H146").Value = 20
$ws.Range("B147").Value = 1341
$ws.Range("C147").Value = 5
$ws.Range("D147").Value = 1092
$ws.Range("E147").Value = 232
$ws.Range("H147").Value = 17
$ws.Range("B150").Value = 1308
$ws.Range("C150").Value = 94
$ws.Range("D150").Value = 136
$ws.Range("E150").Value = 1169
$ws.Range("H150").Value = 3
$ws.Range("B151").Value = 1267
$ws.Range("D151").Value = 1013
$ws.Range("E151").Value = 199
$ws.Range("H151").Value = 55
$ws.Range("B152").Value = 1257
$ws.Range("D152").Value = 788
$ws.Range("E152").Value = 387
$ws.Range("H152").Value = 82
$ws.Range("B157").Value = 1005
$ws.Range("C157").Value = 16
$ws.Range("D157").Value = 869
$ws.Range("E157").Value = 83
$ws.Range("B160").Value = 946
$ws.Range("C160").Value = 43
$ws.Range("E160").Value = 650
$ws.Range("B173").Value = 376
$ws.Range("C173").Value = 1
$ws.Range("D173").Value = 331
$ws.Range("E173").Value = 39
$ws.Range("D213").Value = 13
$ws.Range("H213").Value = 0
$ws.Range("D214").Value = 12
$ws.Range("H214").Value = 1
